$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '27.443.83'
$ws.Range('E2').Value2 = '  -0.71%  '

$ws.Range('D3').Value2 = '1.823.72'
$ws.Range('E3').Value2 = '  -2.06%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value2 = '1.005'
$ws.Range('E4').Value2 = '  -0.61%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value2 = '332.24'
$ws.Range('E5').Value2 = '  -0.48%  '

$ws.Range('E6').Value2 = '  -0.52%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value2 = '0.4581'
$ws.Range('E7').Value2 = '  -2.29%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value2 = '0.3802'
$ws.Range('E8').Value2 = '  -2.79%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value2 = '46.31'
$ws.Range('E9').Value2 = '  +1.18%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value2 = '0.07884'
$ws.Range('E10').Value2 = '  -1.32%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value2 = '0.9689'
$ws.Range('E11').Value2 = '  -3.37%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value2 = '20.99'
$ws.Range('E12').Value2 = '  -3.59%  '

$ws.Range('B13').Value2 = 'Polkadot'
$ws.Range('C13').Value2 = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value2 = '5.879'
$ws.Range('E13').Value2 = '  -1.77%  '

$ws.Range('B14').Value2 = 'WrappedEther'
$ws.Range('C14').Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value2 = '1.803.89'
$ws.Range('E14').Value2 = '  -4.29%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value2 = '7.042'
$ws.Range('E15').Value2 = '  -2.76%  '

$ws.Range('E16').Value2 = '  -0.53%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value2 = '89.52'
$ws.Range('E17').Value2 = '  +1.51%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value2 = '0.06620'
$ws.Range('E18').Value2 = '  -1.74%  '

$ws.Range('E19').Value2 = '  -1.72%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value2 = '17.05'
$ws.Range('E20').Value2 = '  -0.15%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value2 = '1.004'
$ws.Range('E21').Value2 = '  -0.57%  '

$ws.Range('D22').Value2 = '27.440.40'
$ws.Range('E22').Value2 = '  -0.65%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value2 = '5.332'
$ws.Range('E23').Value2 = '  -2.31%  '

$ws.Range('E24').Value2 = '  -0.73%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value2 = '2.306'
$ws.Range('E25').Value2 = '  -0.08%  '

$ws.Range('D26').Value2 = '2.042.05'
$ws.Range('E26').Value2 = '  -2.85%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value2 = '155.58'
$ws.Range('E27').Value2 = '  -2.47%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value2 = '19.35'
$ws.Range('E28').Value2 = '  -2.35%  '

$ws.Range('E29').Value2 = '  -4.37%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value2 = '5.267'
$ws.Range('E30').Value2 = '  -3.00%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value2 = '118.15'
$ws.Range('E31').Value2 = '  -2.73%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value2 = '0.9416'
$ws.Range('E32').Value2 = '  -3.79%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value2 = '0.09298'
$ws.Range('E33').Value2 = '  -2.02%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value2 = '3.583'
$ws.Range('E34').Value2 = '  -0.97%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value2 = '5.244'
$ws.Range('E35').Value2 = '  -0.95%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value2 = '1.316'
$ws.Range('E36').Value2 = '  -1.39%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value2 = '0.05921'
$ws.Range('E37').Value2 = '  -2.36%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value2 = '0.02178'
$ws.Range('E38').Value2 = '  -2.35%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value2 = '8.049'
$ws.Range('E39').Value2 = '  -3.04%  '

$ws.Range('E40').Value2 = '  -3.13%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value2 = '0.5759'
$ws.Range('E41').Value2 = '  -3.15%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value2 = '0.1824'
$ws.Range('E42').Value2 = '  -3.08%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value2 = '9.956'
$ws.Range('E43').Value2 = '  -2.97%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value2 = '1.281'
$ws.Range('E44').Value2 = '  +2.45%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value2 = '0.5438'
$ws.Range('E45').Value2 = '  -3.46%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value2 = '11.93'
$ws.Range('E46').Value2 = '  -2.31%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value2 = '1.862'
$ws.Range('E47').Value2 = '  -3.13%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value2 = '110.87'
$ws.Range('E48').Value2 = '  -1.11%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value2 = '0.06598'
$ws.Range('E49').Value2 = '  -2.46%  '

$ws.Range('E50').Value2 = '  -0.62%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value2 = '1.040'
$ws.Range('E51').Value2 = '  -1.45%  '
